# feat: add 2022-Q4 data
#
# 1. Insert a brand-new worksheet "2022-Q4" right before "2022-Q3" and fill
#    it with the quarterly fund-holding detail rows.
# 2. Insert a new summary row at the top of the "总计" sheet's data (row 2)
#    with the 2022-Q4 totals, pushing the older quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet before the existing "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row
$q4Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q4Headers.Length; $i++) {
    $cell = $q4.Cells.Item(1, $i + 2)
    $cell.Value = $q4Headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Data rows: index, code, name, size, stockPosition, positionRatio, marketValue, rank
$q4Data = @(
    @(0, "720001", "财通价值动量混合", "36.34", "79.52", "4.10", "1.4899", 8),
    @(1, "014915", "财通匠心优选一年持有期混合A", "5.48", "91.42", "4.87", "0.2669", 10),
    @(2, "501046", "财通多策略福鑫定期开放灵活配置混合", "2.69", "91.79", "4.95", "0.1332", 8),
    @(3, "008983", "财通科技创新混合A", "2.87", "87.06", "4.07", "0.1168", 9),
    @(4, "009062", "财通智慧成长混合A", "2.05", "86.49", "4.82", "0.0988", 7),
    @(5, "009063", "财通智慧成长混合C", "1.74", "86.49", "4.82", "0.0839", 7),
    @(6, "015838", "广发招利混合A", "1.42", "90.33", "5.62", "0.0798", 5),
    @(7, "008984", "财通科技创新混合C", "1.79", "87.06", "4.07", "0.0729", 9),
    @(8, "015839", "广发招利混合C", "0.58", "90.33", "5.62", "0.0326", 5),
    @(9, "014916", "财通匠心优选一年持有期混合C", "0.59", "91.42", "4.87", "0.0287", 10),
    @(10, "008135", "华宸未来价值先锋混合", "0.30", "83.80", "4.53", "0.0136", 10)
)

foreach ($row in $q4Data) {
    $r = [int]$row[0] + 2

    $idxCell = $q4.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $bCell = $q4.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $q4.Cells.Item($r, 3).Value = $row[2]

    # Columns D, E, F, G hold text (keep trailing zeros / decimal formatting
    # intact) - force text format before assignment so Excel doesn't coerce
    # the numeric-looking strings into real numbers.
    $dCell = $q4.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $eCell = $q4.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]

    $fCell = $q4.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[5]

    $gCell = $q4.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[6]

    $q4.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 summary row into "总计"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

$summaryData = @(
    @(0, "2022-Q4", 11, 2.42),
    @(1, "2022-Q3", 6, 0.75),
    @(2, "2022-Q2", 2, 0.11),
    @(3, "2021-Q3", 3, 0.24),
    @(4, "2021-Q2", 2, 0)
)

foreach ($row in $summaryData) {
    $r = [int]$row[0] + 2

    $idxCell = $summary.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# Restore the originally-active sheet/tab (adding the new sheet made itself
# active as a side effect).
$wb.Worksheets.Item("2021-Q2").Activate()
